$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)
$new = $wb.Worksheets.Add($firstSheet)
$new.Name = "Browser"
$new.Range("A1").Value = "BrowserName"
$new.Range("A2").Value = "chrome"
$null = $new.Range("A2").Select()
